$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before current row 26 (old row26 "支票用途" becomes row28).
$ws.Range("A26:G27").Insert()

# Copy formatting (borders/fonts/number formats) from row 25, a plain
# unmerged data row with no special row height, onto the two new rows.
# (Columns A:E and G are copied separately so an incidental blank F cell
# isn't materialized - row 25 never had an F cell to begin with.)
$ws.Range("A25:E25").Copy()
$ws.Range("A26:E27").PasteSpecial(-4122)
$ws.Range("G25").Copy()
$ws.Range("G26:G27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new unique shared strings in the same order the author
# entered them: 支票銀行, 支票分行, NVARCHAR2, BankItem, BranchItem.
$ws.Range("C26").Value = "支票銀行"
$ws.Range("C27").Value = "支票分行"
$ws.Range("D26").Value = "NVARCHAR2"
$ws.Range("D27").Value = "NVARCHAR2"
$ws.Range("B26").Value = "BankItem"
$ws.Range("B27").Value = "BranchItem"

# Row 26: BankItem / 支票銀行 / NVARCHAR2 / 50
$ws.Range("E26").Value = 50
$ws.Range("G26").Value = ""

# Row 27: BranchItem / 支票分行 / NVARCHAR2 / 50
$ws.Range("E27").Value = 50
$ws.Range("G27").Value = ""

# Re-enter the running "+1" counter formula for the new rows and every row
# pushed down below them, one cell at a time, so every cell keeps a valid,
# self-contained formula (and correct cached value) after the insert.
for ($r = 26; $r -le 42; $r++) {
  $prev = $r - 1
  $ws.Range("A$r").Formula = "=A$prev+1"
}

# Update the view: scroll/selection moved while editing this part of the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D25").Select()
